$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / URL / percentage cells - direct assignment
$ws.Range('D2').Value = '28.045.14'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '1.888.76'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  -4.48%  '
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('E9').Value = '  -5.45%  '
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '1.883.79'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('E15').Value = '  -3.48%  '
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').Value = '28.110.73'
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.104.03'
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E31').Value = '  -3.41%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E35').Value = '  -5.19%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E36').Value = '  -3.05%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('E38').Value = '  +11.00%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E40').Value = '  -4.26%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E51').Value = '  -3.63%  '

# Numeric-looking Price cells must stay as TEXT (match source formatting,
# e.g. leading/trailing zeros, exact decimal digits). Force text format,
# assign, then restore default "Normal" style so no stray style persists.
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '313.98'
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.5004'
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3875'
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.09172'
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '1.126'
$c.Style = "Normal"
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '41.72'
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '6.353'
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.77'
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '7.300'
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '91.61'
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.06633'
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '17.86'
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.227'
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '11.36'
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.315'
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.549'
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '158.43'
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '20.74'
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '126.81'
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.072'
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.1051'
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.586'
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '3.601'
$c.Style = "Normal"
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '9.392'
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.06560'
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.02406'
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.322'
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.2186'
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.213'
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.6414'
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '11.54'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '4.941'
$c.Style = "Normal"
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '13.39'
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.6025'
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.301'
$c.Style = "Normal"
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '3.682'
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '1.995'
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.205'
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '120.73'
$c.Style = "Normal"
